# Update cryptocurrency price/volume snapshot with latest scraped values.
# Generated from GitHub Actions run: "Updated symbol list on Thu Jan 26 19:53:26 UTC 2023 with GitHub Actions"
#
# The Price (column D) and Volume(1h) (column E) cells are stored as text
# (e.g. "305.09", "1.06%"), so a leading apostrophe is used to force Excel
# to keep the exact text representation instead of re-interpreting the
# value as a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'305.09"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("D3").Value = "'35.92"
$ws.Range("E3").Value = "'0.99%"
$ws.Range("D4").Value = "'5.026"
$ws.Range("E4").Value = "'-1.06%"
$ws.Range("E5").Value = "'0.87%"
$ws.Range("D6").Value = "'1.950"
$ws.Range("E6").Value = "'0.55%"
$ws.Range("D7").Value = "'4.137"
$ws.Range("E7").Value = "'2.19%"
$ws.Range("D8").Value = "'7.846"
$ws.Range("E8").Value = "'0.90%"
$ws.Range("D9").Value = "'0.9307"
$ws.Range("E9").Value = "'0.23%"
$ws.Range("D10").Value = "'0.1257"
$ws.Range("E10").Value = "'-21.40%"
$ws.Range("D11").Value = "'0.1910"
$ws.Range("E11").Value = "'0.48%"
$ws.Range("D12").Value = "'0.09228"
$ws.Range("E12").Value = "'2.92%"
$ws.Range("D13").Value = "'0.03509"
$ws.Range("E13").Value = "'1.63%"
$ws.Range("E14").Value = "'0.36%"
$ws.Range("D15").Value = "'0.001423"
$ws.Range("E15").Value = "'1.91%"
$ws.Range("D16").Value = "'0.006694"
$ws.Range("E16").Value = "'16.64%"
$ws.Range("D17").Value = "'3.615"
$ws.Range("E17").Value = "'2.17%"
$ws.Range("D18").Value = "'3.085"
$ws.Range("E18").Value = "'7.11%"
$ws.Range("E19").Value = "'-0.04%"
$ws.Range("D20").Value = "'5.168"
$ws.Range("E20").Value = "'2.36%"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("E21").Value = "'-1.12%"
$ws.Range("E22").Value = "'5.56%"
$ws.Range("E23").Value = "'-2.18%"
$ws.Range("E24").Value = "'1.91%"
$ws.Range("D25").Value = "'0.004727"
$ws.Range("E25").Value = "'-1.17%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'5.80%"
$ws.Range("D27").Value = "'0.0003129"
$ws.Range("E27").Value = "'3.53%"
$ws.Range("D39").Value = "'0.01969"
$ws.Range("E39").Value = "'6.32%"
$ws.Range("D40").Value = "'0.05163"
$ws.Range("E40").Value = "'7.68%"
$ws.Range("D41").Value = "'0.007564"
$ws.Range("E41").Value = "'3.42%"
$ws.Range("D42").Value = "'0.01012"
$ws.Range("E42").Value = "'-4.78%"
$ws.Range("D43").Value = "'0.1373"
$ws.Range("E43").Value = "'3.13%"
$ws.Range("E44").Value = "'-0.37%"
$ws.Range("D45").Value = "'0.01065"
$ws.Range("E45").Value = "'9.84%"
$ws.Range("D46").Value = "'0.00006391"
$ws.Range("E46").Value = "'2.65%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("D48").Value = "'64.96"
$ws.Range("E48").Value = "'0.45%"
$ws.Range("D49").Value = "'0.001600"
$ws.Range("E49").Value = "'-3.55%"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("E51").Value = "'0.11%"
